$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new journal entry as row 11, matching the formatting of row 10
# (date column style + wrapped-text description column style).
$ws.Range("A10:B10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A11").Value = 44993
$ws.Range("B11").Value = "rencontre avec M Hurni pour une explication du fonctionement des Models (mon sauveur)"

$ws.Rows.Item(11).RowHeight = $ws.Rows.Item(10).RowHeight

$ws.Range("B11").Select()
